$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column A (test_number) values from 3 to 4 for rows 2-11
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = 4
}

# Update the active selection to match the recorded cursor position
$ws.Range("H11").Select()
